$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 3542.5715
$ws.Range("J62").Value = 3501.25
$ws.Range("L62").Value = 3501.25
$ws.Range("N62").Value = -4749.25
# Row 65
$ws.Range("H65").Value = 3542.5715
$ws.Range("J65").Value = 3501.25
$ws.Range("L65").Value = 17506.25
$ws.Range("N65").Value = -23746.25
# Row 92
$ws.Range("H92").Value = 1151.0769
$ws.Range("I92").Value = 262.375
$ws.Range("J92").Value = 2573
$ws.Range("K92").Value = 262.375
$ws.Range("L92").Value = 2573
$ws.Range("M92").Value = 985.625
$ws.Range("N92").Value = -5069
# Row 94
$ws.Range("H94").Value = 3597
$ws.Range("I94").Value = 3597
$ws.Range("K94").Value = 3597
$ws.Range("M94").Value = -3146
# Row 99
$ws.Range("H99").Value = 1015.75
$ws.Range("I99").Value = 288
$ws.Range("J99").Value = 1743.5
$ws.Range("K99").Value = 864
$ws.Range("L99").Value = 5230.5
$ws.Range("M99").Value = 634
$ws.Range("N99").Value = -8226.5
# Row 100
$ws.Range("H100").Value = 1403
$ws.Range("I100").Value = 1367.3334
$ws.Range("J100").Value = 1438.6666
$ws.Range("K100").Value = 1367.3334
$ws.Range("L100").Value = 1438.6666
$ws.Range("M100").Value = -826.3334
$ws.Range("N100").Value = -2520.6666
# Row 101
$ws.Range("H101").Value = 1120.3077
$ws.Range("I101").Value = 442
$ws.Range("J101").Value = 1911.6666
$ws.Range("K101").Value = 1326
$ws.Range("L101").Value = 5734.9998
$ws.Range("M101").Value = 296
$ws.Range("N101").Value = -8978.9998
# Row 103
$ws.Range("H103").Value = 1145.5769
$ws.Range("I103").Value = 766.44446
$ws.Range("J103").Value = 1346.2941
$ws.Range("K103").Value = 2299.33338
$ws.Range("L103").Value = 4038.8823
$ws.Range("M103").Value = -1713.33338
$ws.Range("N103").Value = -5210.8823
# Row 106
$ws.Range("H106").Value = 2600.8333
$ws.Range("I106").Value = 2601.25
$ws.Range("J106").Value = 2600
$ws.Range("K106").Value = 2601.25
$ws.Range("L106").Value = 2600
$ws.Range("M106").Value = -1970.25
$ws.Range("N106").Value = -3862
# Row 121
$ws.Range("H121").Value = 3529.3333
$ws.Range("J121").Value = 3529.3333
$ws.Range("L121").Value = 10587.9999
$ws.Range("N121").Value = -14081.9999
# Row 125
$ws.Range("H125").Value = 3063.3333
$ws.Range("I125").Value = 1178
$ws.Range("J125").Value = 4006
$ws.Range("K125").Value = 10602
$ws.Range("L125").Value = 36054
$ws.Range("M125").Value = -8142
$ws.Range("N125").Value = -40974
# Row 128
$ws.Range("H128").Value = 90000
$ws.Range("J128").Value = 90000
$ws.Range("L128").Value = 90000
$ws.Range("N128").Value = -99960

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 2306.878
$ws.Range("I61").Value = 1747.6207
$ws.Range("K61").Value = 1747.6207
$ws.Range("M61").Value = -1535.6207
# Row 74
$ws.Range("H74").Value = 5067.4
$ws.Range("I74").Value = 5824.25
$ws.Range("J74").Value = 2040
$ws.Range("K74").Value = 5824.25
$ws.Range("L74").Value = 2040
$ws.Range("M74").Value = -4950.25
$ws.Range("N74").Value = -3788
# Row 77
$ws.Range("H77").Value = 5067.4
$ws.Range("I77").Value = 5824.25
$ws.Range("J77").Value = 2040
$ws.Range("K77").Value = 29121.25
$ws.Range("L77").Value = 10200
$ws.Range("M77").Value = -24753.25
$ws.Range("N77").Value = -18936
# Row 97
$ws.Range("H97").Value = 616.9375
$ws.Range("I97").Value = 489.92307
$ws.Range("J97").Value = 1167.3334
$ws.Range("K97").Value = 489.92307
$ws.Range("L97").Value = 1167.3334
$ws.Range("M97").Value = 6.076930000000004
$ws.Range("N97").Value = -2159.3334
# Row 102
$ws.Range("H102").Value = 3128.5715
$ws.Range("I102").Value = 2824
$ws.Range("K102").Value = 2824
$ws.Range("M102").Value = -1202
# Row 122
$ws.Range("H122").Value = 2605773.8
$ws.Range("I122").Value = 3126528.5
$ws.Range("J122").Value = 2000
$ws.Range("K122").Value = 9379585.5
$ws.Range("L122").Value = 6000
$ws.Range("M122").Value = -9377135.5
$ws.Range("N122").Value = -10900
# Row 124
$ws.Range("H124").Value = 16600
$ws.Range("J124").Value = 16600
$ws.Range("L124").Value = 16600
$ws.Range("N124").Value = -26420
# Row 136
$ws.Range("H136").Value = 2306.878
$ws.Range("I136").Value = 1747.6207
$ws.Range("K136").Value = 5242.8621
$ws.Range("M136").Value = -2692.8621
# Row 138
$ws.Range("H138").Value = 30539.666
$ws.Range("J138").Value = 30539.666
$ws.Range("L138").Value = 30539.666
$ws.Range("N138").Value = -40819.666

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 692.5714
$ws.Range("I94").Value = 634.1111
$ws.Range("J94").Value = 1043.3334
$ws.Range("K94").Value = 634.1111
$ws.Range("L94").Value = 1043.3334
$ws.Range("M94").Value = -183.1111
$ws.Range("N94").Value = -1945.3334
# Row 134
$ws.Range("H134").Value = 2410.1667
$ws.Range("I134").Value = 2368.3333
$ws.Range("J134").Value = 2535.6667
$ws.Range("K134").Value = 7104.999899999999
$ws.Range("L134").Value = 7607.000100000001
$ws.Range("M134").Value = -4569.999899999999
$ws.Range("N134").Value = -12677.0001

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 4977.1
$ws.Range("I31").Value = 5106.75
$ws.Range("J31").Value = 4828.9287
$ws.Range("K31").Value = 5106.75
$ws.Range("L31").Value = 4828.9287
$ws.Range("M31").Value = -4811.75
$ws.Range("N31").Value = -5418.9287
# Row 34
$ws.Range("H34").Value = 4977.1
$ws.Range("I34").Value = 5106.75
$ws.Range("J34").Value = 4828.9287
$ws.Range("K34").Value = 5106.75
$ws.Range("L34").Value = 4828.9287
$ws.Range("M34").Value = -4904.75
$ws.Range("N34").Value = -5232.9287
# Row 105
$ws.Range("H105").Value = 51798.332
$ws.Range("I105").Value = 1811.3334
$ws.Range("J105").Value = 301733.34
$ws.Range("K105").Value = 1811.3334
$ws.Range("L105").Value = 301733.34
$ws.Range("M105").Value = -64.33339999999998
$ws.Range("N105").Value = -305227.34
# Row 135
$ws.Range("H135").Value = 41259.832
$ws.Range("J135").Value = 41259.832
$ws.Range("L135").Value = 41259.832
$ws.Range("N135").Value = -51399.832

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 122
$ws.Range("H122").Value = 799.8461
$ws.Range("J122").Value = 1149.75
$ws.Range("L122").Value = 10347.75
$ws.Range("N122").Value = -15247.75
# Row 131
$ws.Range("H131").Value = 721.6393399999999
$ws.Range("J131").Value = 962.6579
$ws.Range("L131").Value = 2887.9737
$ws.Range("N131").Value = -12967.9737
# Row 134
$ws.Range("H134").Value = 3821.5789
$ws.Range("I134").Value = 1861
$ws.Range("K134").Value = 5583
$ws.Range("M134").Value = -513

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 97
$ws.Range("H97").Value = 500884.34
$ws.Range("I97").Value = 527194.0600000001
$ws.Range("J97").Value = 1000
$ws.Range("K97").Value = 527194.0600000001
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = -526698.0600000001
$ws.Range("N97").Value = -1992
# Row 113
$ws.Range("H113").Value = 12156.909
$ws.Range("I113").Value = 1900
$ws.Range("J113").Value = 24465.2
$ws.Range("K113").Value = 1900
$ws.Range("L113").Value = 24465.2
$ws.Range("M113").Value = 270
$ws.Range("N113").Value = -28805.2
# Row 122
$ws.Range("H122").Value = 33334300
$ws.Range("I122").Value = 33334300
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 100002900
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -100000450
$ws.Range("N122").ClearContents()
# Row 123
$ws.Range("H123").Value = 18604.525
$ws.Range("J123").Value = 18604.525
$ws.Range("L123").Value = 18604.525
$ws.Range("N123").Value = -23504.525

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 70
$ws.Range("H70").Value = 10000
$ws.Range("J70").Value = 10000
$ws.Range("L70").Value = 10000
$ws.Range("N70").Value = -10540
# Row 73
$ws.Range("H73").Value = 10000
$ws.Range("J73").Value = 10000
$ws.Range("L73").Value = 10000
$ws.Range("N73").Value = -11872
# Row 93
$ws.Range("H93").Value = 2722.5
$ws.Range("I93").Value = 2722.5
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 2722.5
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -1474.5
$ws.Range("N93").ClearContents()
# Row 125
$ws.Range("H125").Value = 63216.43
$ws.Range("J125").Value = 63216.43
$ws.Range("L125").Value = 63216.43
$ws.Range("N125").Value = -73056.42999999999
# Row 136
$ws.Range("H136").Value = 1839.6
$ws.Range("I136").Value = 1466.1666
$ws.Range("J136").Value = 3333.3333
$ws.Range("K136").Value = 4398.4998
$ws.Range("L136").Value = 9999.999899999999
$ws.Range("M136").Value = -1848.4998
$ws.Range("N136").Value = -15099.9999

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 27
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("N27").ClearContents()
# Row 96
$ws.Range("H96").Value = 2923.8572
$ws.Range("I96").Value = 1703
$ws.Range("J96").Value = 3127.3333
$ws.Range("K96").Value = 1703
$ws.Range("L96").Value = 3127.3333
$ws.Range("M96").Value = -330
$ws.Range("N96").Value = -5873.3333
# Row 100
$ws.Range("H100").Value = 2195
$ws.Range("I100").Value = 2801
$ws.Range("J100").Value = 1952.6
$ws.Range("K100").Value = 5602
$ws.Range("L100").Value = 3905.2
$ws.Range("M100").Value = -5061
$ws.Range("N100").Value = -4987.2
# Row 115
$ws.Range("H115").Value = 38999
$ws.Range("J115").Value = 38999
$ws.Range("L115").Value = 38999
$ws.Range("N115").Value = -42133
# Row 140
$ws.Range("H140").Value = 27067
$ws.Range("J140").Value = 27067
$ws.Range("L140").Value = 27067
$ws.Range("N140").Value = -37427
